$wb = $excel.ActiveWorkbook

# ALC!row11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 96.59999999999999
$ws.Range("I11").Value = 96.59999999999999
$ws.Range("K11").Value = 96.59999999999999
$ws.Range("M11").Value = 43.40000000000001

# ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8485.583000000001
$ws.Range("I40").Value = 6466.2
$ws.Range("J40").Value = 9928
$ws.Range("K40").Value = 6466.2
$ws.Range("L40").Value = 9928
$ws.Range("M40").Value = -6291.2
$ws.Range("N40").Value = -10278

# ALC!row88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1674
$ws.Range("I88").Value = 1674
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1674
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1268
$ws.Range("N88").ClearContents()

# ALC!row91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1674
$ws.Range("I91").Value = 1674
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1674
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -270
$ws.Range("N91").ClearContents()

# ALC!row96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 177.90909
$ws.Range("I96").Value = 226.71428
$ws.Range("J96").Value = 92.5
$ws.Range("K96").Value = 680.14284
$ws.Range("L96").Value = 277.5
$ws.Range("M96").Value = 692.85716
$ws.Range("N96").Value = -3023.5

# ALC!row97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 707.5
$ws.Range("J97").Value = 594.3333
$ws.Range("L97").Value = 1782.9999
$ws.Range("N97").Value = -2774.9999

# ALC!row99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 6506.75
$ws.Range("I99").Value = 1164
$ws.Range("K99").Value = 3492
$ws.Range("M99").Value = -1994

# ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1185.8462
$ws.Range("I111").Value = 418
$ws.Range("J111").Value = 1527.1111
$ws.Range("K111").Value = 1254
$ws.Range("L111").Value = 4581.3333
$ws.Range("M111").Value = 1813
$ws.Range("N111").Value = -10715.3333

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 9029.343000000001
$ws.Range("I132").Value = 7863
$ws.Range("K132").Value = 23589
$ws.Range("M132").Value = -21059

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4121.35
$ws.Range("I137").Value = 905.5714
$ws.Range("K137").Value = 2716.7142
$ws.Range("M137").Value = -166.7142000000003

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1133.9412
$ws.Range("I97").Value = 1127.3572
$ws.Range("K97").Value = 1127.3572
$ws.Range("M97").Value = -631.3571999999999

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1603.4
$ws.Range("I132").Value = 1821.75
$ws.Range("K132").Value = 5465.25
$ws.Range("M132").Value = -2935.25

# CRP!row3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3333467.8
$ws.Range("J3").Value = 201
$ws.Range("L3").Value = 201
$ws.Range("N3").Value = -427

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 987.53845
$ws.Range("I122").Value = 987.53845
$ws.Range("K122").Value = 2962.61535
$ws.Range("M122").Value = -512.61535

# CUL!row14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 701.75
$ws.Range("I14").Value = 701.75
$ws.Range("K14").Value = 2105.25
$ws.Range("M14").Value = -1932.25

# CUL!row38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 445.4737
$ws.Range("I38").Value = 447.75
$ws.Range("K38").Value = 1343.25
$ws.Range("M38").Value = -996.25

# CUL!row127
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# GSM!row5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5402.5
$ws.Range("I80").Value = 5402.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5402.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -4404.5
$ws.Range("N80").ClearContents()

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5402.5
$ws.Range("I83").Value = 5402.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 27012.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -22020.5
$ws.Range("N83").ClearContents()

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2266.6667
$ws.Range("I16").Value = 2266.6667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2266.6667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2096.6667
$ws.Range("N16").ClearContents()

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7584.3335
$ws.Range("J40").Value = 10169
$ws.Range("L40").Value = 10169
$ws.Range("N40").Value = -10441

# LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6265.8335
$ws.Range("I61").Value = 4532
$ws.Range("J61").Value = 7999.6665
$ws.Range("K61").Value = 4532
$ws.Range("L61").Value = 7999.6665
$ws.Range("M61").Value = -4330
$ws.Range("N61").Value = -8403.666499999999

# LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6265.8335
$ws.Range("I113").Value = 4532
$ws.Range("J113").Value = 7999.6665
$ws.Range("K113").Value = 4532
$ws.Range("L113").Value = 7999.6665
$ws.Range("M113").Value = -2362
$ws.Range("N113").Value = -12339.6665

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2814.1667
$ws.Range("I136").Value = 2377
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7131
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4581
$ws.Range("N136").Value = -20100

# WVR!row62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9507.916999999999
$ws.Range("I62").Value = 7199.25
$ws.Range("J62").Value = 10662.25
$ws.Range("K62").Value = 7199.25
$ws.Range("L62").Value = 10662.25
$ws.Range("M62").Value = -6575.25
$ws.Range("N62").Value = -11910.25

# WVR!row65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9507.916999999999
$ws.Range("I65").Value = 7199.25
$ws.Range("J65").Value = 10662.25
$ws.Range("K65").Value = 35996.25
$ws.Range("L65").Value = 53311.25
$ws.Range("M65").Value = -32876.25
$ws.Range("N65").Value = -59551.25

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 443.25
$ws.Range("I81").Value = 445
$ws.Range("J81").Value = 441.5
$ws.Range("K81").Value = 890
$ws.Range("L81").Value = 883
$ws.Range("M81").Value = 171
$ws.Range("N81").Value = -3005

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 443.25
$ws.Range("I84").Value = 445
$ws.Range("J84").Value = 441.5
$ws.Range("K84").Value = 4450
$ws.Range("L84").Value = 4415
$ws.Range("M84").Value = 854
$ws.Range("N84").Value = -15023

# WVR!row125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 28999.5
$ws.Range("J125").Value = 28999.5
$ws.Range("L125").Value = 28999.5
$ws.Range("N125").Value = -38839.5
